# Updates the EC (Estado de Cuenta) worksheet:
#  - Adds a new "Periodo Mora" row (2508) as part 1 of new estado de cuenta
#  - Updates existing period rows / totals
#  - Shifts the signature block down one row
#  - Nudges the logo image slightly to the left

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new data row (19) below the current last period row (18).
#    Copying the whole row brings along the "last row" border styling so the
#    new row keeps the closing bottom border of the table.
# ---------------------------------------------------------------------------
$ws.Rows.Item(18).Copy()
$ws.Rows.Item(19).Insert()

# Row 18 is no longer the last row of the table, so it should take on the
# "interior" row styling that rows 16/17 already use.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122) # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Update the period values for rows 16-19.
# ---------------------------------------------------------------------------
$ws.Range("E16").Value = "2505"
$ws.Range("F16").Value = 25600

$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 64000

$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 64000

$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 64000
$ws.Range("G19").Value = 1600000

# ---------------------------------------------------------------------------
# 3) Update totals: Valor Mora (E11) and Cant. Periodos (F13).
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 217600
$ws.Range("F13").Value = 4

# ---------------------------------------------------------------------------
# 4) Nudge the logo image left by 171450 EMU (~0.476 cm).
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = $shp.Left - 4.5
